$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# The title currently holds two runs: "Elastic Search " (no dirty attr)
# and "Essentials" (dirty="0"). Remove the first run's characters so only
# the second run (with its "dirty=0" formatting) remains, then set the
# full text on the now-single run so it keeps that run's formatting.
$firstRunLength = "Elastic Search ".Length
$c1 = $tr.Characters(1, $firstRunLength)
$c1.Text = ""

$tr2 = $sh.TextFrame.TextRange
$tr2.Text = "Elastic Search Essentials"
